# fix generation fichier excel
# Rebuild the worksheet data: replace the old 5-row demo dataset with the
# new 16-row dataset (dates/times are written as text, the "numero" column
# mixes true numbers and text codes exactly as the source data does).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("admin", "16/05/2012", "12:30:08", 12345,    "DELL300"),
    @("admin", "16/05/2012", "12:35:29", "AZ123456","DELL200"),
    @("admin", "16/05/2012", "13:24:02", 123456,   "DELL400"),
    @("admin", "16/05/2012", "13:24:24", "AZ123456","BOUY300"),
    @("admin", "16/05/2012", "13:24:43", 12345,    "DELL300"),
    @("admin", "16/05/2012", "13:26:26", 12345,    "DELL300"),
    @("admin", "16/05/2012", "13:44:51", "AZERTY", "DELL100"),
    @("admin", "16/05/2012", "16:43:40", 12345,    "DELL300"),
    @("admin", "16/05/2012", "17:19:40", 12345,    "DELL300"),
    @("admin", "16/05/2012", "17:19:40", "AZ123456","BOUY300"),
    @("admin", "16/05/2012", "17:19:40", 12345,    "DELL700"),
    @("admin", "16/05/2012", "17:19:40", "AZ123456","BOUY700"),
    @("admin", "16/05/2012", "17:31:16", 12345,    "BOUY800"),
    @("admin", "16/05/2012", "17:31:16", "AZ123456","BOUY800"),
    @("admin", "16/05/2012", "17:31:16", 12345,    "DELL700"),
    @("admin", "16/05/2012", "17:31:16", "AZ123456","BOUY700")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 1
    $values = $data[$i]
    for ($c = 0; $c -lt $values.Length; $c++) {
        $col = $c + 1
        $ws.Cells.Item($row, $col).Value = $values[$c]
    }
}
